# "added albums loading and post history"
#
# Mark the checklist rows related to albums/post-history functionality
# as completed by filling "YES" into column B for the relevant rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$completedRows = @(41, 42, 43, 71, 87, 88, 91, 92, 95)
foreach ($r in $completedRows) {
    $ws.Cells.Item($r, 2).Value = "YES"
}

# Leave the view scrolled/selected where the author last worked.
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D41").Select()
